$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('E2').Value = '2026-02-08 16:48:24'
$ws.Range('I2').Value = '2.3 mm'
$ws.Range('E3').Value = '2026-02-08 16:48:26'
$ws.Range('I3').Value = '1.6 mm'
$ws.Range('L3').Value = '49.7 km/h - 236º 16:28 TU'
$ws.Range('E4').Value = '2026-02-08 16:48:29'
$ws.Range('K4').Value = '9.4 MJ/m2'
$ws.Range('E5').Value = '2026-02-08 16:48:31'
$ws.Range('H5').NumberFormat = '@'
$ws.Range('H5').Value = '89%'
$ws.Range('I5').Value = '4.4 mm'
$ws.Range('K5').Value = '5.4 MJ/m2'
$ws.Range('E6').Value = '2026-02-08 16:48:34'
$ws.Range('K6').Value = '12.1 MJ/m2'
$ws.Range('E7').Value = '2026-02-08 16:48:36'
$ws.Range('K7').Value = '11.2 MJ/m2'
$ws.Range('O7').Value = '12.7 °C'
$ws.Range('E8').Value = '2026-02-08 16:48:38'
$ws.Range('H8').NumberFormat = '@'
$ws.Range('H8').Value = '77%'
$ws.Range('E9').Value = '2026-02-08 16:48:41'
$ws.Range('K9').Value = '12.1 MJ/m2'
$ws.Range('O9').Value = '10.4 °C'
$ws.Range('E10').Value = '2026-02-08 16:48:43'
$ws.Range('H10').NumberFormat = '@'
$ws.Range('H10').Value = '81%'
$ws.Range('E11').Value = '2026-02-08 16:48:45'
$ws.Range('H11').NumberFormat = '@'
$ws.Range('H11').Value = '80%'
$ws.Range('O11').Value = '4.4 °C'
$ws.Range('E12').Value = '2026-02-08 16:48:48'
$ws.Range('E13').Value = '2026-02-08 16:48:50'
$ws.Range('K13').Value = '10.8 MJ/m2'
$ws.Range('O13').Value = '3.2 °C'
$ws.Range('E14').Value = '2026-02-08 16:48:52'
$ws.Range('H14').NumberFormat = '@'
$ws.Range('H14').Value = '78%'
$ws.Range('K14').Value = '11.8 MJ/m2'
$ws.Range('O14').Value = '10.9 °C'
$ws.Range('E15').Value = '2026-02-08 16:48:55'
$ws.Range('O15').Value = '10.0 °C'
$ws.Range('E16').Value = '2026-02-08 16:48:57'
$ws.Range('I16').Value = '1.7 mm'
$ws.Range('K16').Value = '7.5 MJ/m2'
$ws.Range('L16').Value = '60.8 km/h - 231º 16:12 TU'
$ws.Range('E17').Value = '2026-02-08 16:48:59'
$ws.Range('K17').Value = '7.6 MJ/m2'
$ws.Range('O17').Value = '-0.1 °C'
$ws.Range('E18').Value = '2026-02-08 16:49:02'
$ws.Range('J18').Value = '1001.6 hPa'
$ws.Range('K18').Value = '11.1 MJ/m2'
$ws.Range('E19').Value = '2026-02-08 16:49:04'
$ws.Range('I19').Value = '8.0 mm'
$ws.Range('K19').Value = '10.2 MJ/m2'
$ws.Range('O19').Value = '4.7 °C'
$ws.Range('E20').Value = '2026-02-08 16:49:06'
$ws.Range('K20').Value = '10.1 MJ/m2'
$ws.Range('E21').Value = '2026-02-08 16:49:09'
$ws.Range('H21').NumberFormat = '@'
$ws.Range('H21').Value = '83%'
$ws.Range('J21').Value = '1002.7 hPa'
$ws.Range('K21').Value = '11.5 MJ/m2'
$ws.Range('L21').Value = '22.3 km/h - 354º 16:07 TU'
$ws.Range('O21').Value = '5.1 °C'
$ws.Range('E22').Value = '2026-02-08 16:49:11'
$ws.Range('K22').Value = '7.4 MJ/m2'
$ws.Range('L22').Value = '45.0 km/h - 339º 16:23 TU'
$ws.Range('E23').Value = '2026-02-08 16:49:14'
$ws.Range('I23').Value = '3.3 mm'
$ws.Range('K23').Value = '10.3 MJ/m2'
$ws.Range('L23').Value = '43.9 km/h - 328º 16:22 TU'
$ws.Range('E24').Value = '2026-02-08 16:49:16'
$ws.Range('H24').NumberFormat = '@'
$ws.Range('H24').Value = '84%'
$ws.Range('K24').Value = '11.2 MJ/m2'
$ws.Range('O24').Value = '8.2 °C'
$ws.Range('E25').Value = '2026-02-08 16:49:19'
$ws.Range('K25').Value = '9.2 MJ/m2'
$ws.Range('O25').Value = '-2.9 °C'
$ws.Range('E26').Value = '2026-02-08 16:49:21'
$ws.Range('K26').Value = '11.8 MJ/m2'
$ws.Range('E27').Value = '2026-02-08 16:49:24'
$ws.Range('K27').Value = '7.7 MJ/m2'
$ws.Range('E28').Value = '2026-02-08 16:49:26'
$ws.Range('J28').Value = '1001.3 hPa'
$ws.Range('O28').Value = '8.5 °C'
$ws.Range('E29').Value = '2026-02-08 16:49:28'
$ws.Range('H29').NumberFormat = '@'
$ws.Range('H29').Value = '82%'
$ws.Range('K29').Value = '11.7 MJ/m2'
$ws.Range('O29').Value = '10.6 °C'
$ws.Range('E30').Value = '2026-02-08 16:49:31'
$ws.Range('J30').Value = '1001.6 hPa'
$ws.Range('K30').Value = '10.3 MJ/m2'
$ws.Range('E31').Value = '2026-02-08 16:49:33'
$ws.Range('H31').NumberFormat = '@'
$ws.Range('H31').Value = '72%'
$ws.Range('J31').Value = '1000.6 hPa'
$ws.Range('K31').Value = '9.1 MJ/m2'
$ws.Range('N31').Value = '8.9 °C 16:21 TU'
$ws.Range('E32').Value = '2026-02-08 16:49:35'
$ws.Range('K32').Value = '5.1 MJ/m2'
$ws.Range('E33').Value = '2026-02-08 16:49:37'
$ws.Range('K33').Value = '8.9 MJ/m2'
$ws.Range('O33').Value = '2.6 °C'
$ws.Range('E34').Value = '2026-02-08 16:49:40'
$ws.Range('K34').Value = '12.5 MJ/m2'
$ws.Range('O34').Value = '-0.6 °C'
$ws.Range('E35').Value = '2026-02-08 16:49:42'
$ws.Range('K35').Value = '4.6 MJ/m2'
$ws.Range('O35').Value = '3.8 °C'
$ws.Range('E36').Value = '2026-02-08 16:49:45'
$ws.Range('K36').Value = '11.4 MJ/m2'
$ws.Range('E37').Value = '2026-02-08 16:49:47'
$ws.Range('J37').Value = '1002.6 hPa'
$ws.Range('E38').Value = '2026-02-08 16:49:50'
$ws.Range('H38').NumberFormat = '@'
$ws.Range('H38').Value = '74%'
$ws.Range('K38').Value = '10.0 MJ/m2'
$ws.Range('O38').Value = '9.6 °C'
$ws.Range('E39').Value = '2026-02-08 16:49:52'
$ws.Range('K39').Value = '12.7 MJ/m2'
$ws.Range('E40').Value = '2026-02-08 16:49:54'
$ws.Range('H40').NumberFormat = '@'
$ws.Range('H40').Value = '85%'
$ws.Range('O40').Value = '5.6 °C'
$ws.Range('E41').Value = '2026-02-08 16:49:56'
$ws.Range('H41').NumberFormat = '@'
$ws.Range('H41').Value = '73%'
$ws.Range('J41').Value = '1001.8 hPa'
$ws.Range('K41').Value = '12.5 MJ/m2'
$ws.Range('O41').Value = '11.7 °C'
$ws.Range('E42').Value = '2026-02-08 16:49:59'
$ws.Range('E43').Value = '2026-02-08 16:50:01'
$ws.Range('K43').Value = '11.0 MJ/m2'
$ws.Range('O43').Value = '6.5 °C'
$ws.Range('E44').Value = '2026-02-08 16:50:03'
$ws.Range('I44').Value = '1.4 mm'
$ws.Range('K44').Value = '6.0 MJ/m2'
$ws.Range('E45').Value = '2026-02-08 16:50:05'
$ws.Range('G45').Value = '2 cm'
$ws.Range('H45').NumberFormat = '@'
$ws.Range('H45').Value = '74%'
$ws.Range('I45').Value = '1.2 mm'
$ws.Range('J45').Value = '1003.7 hPa'
$ws.Range('E46').Value = '2026-02-08 16:50:08'
$ws.Range('H46').NumberFormat = '@'
$ws.Range('H46').Value = '77%'
$ws.Range('J46').Value = '1003.4 hPa'
$ws.Range('K46').Value = '7.9 MJ/m2'
$ws.Range('O46').Value = '9.0 °C'
